$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4585.3
$ws.Range("I74").Value = 4436.5
$ws.Range("K74").Value = 4436.5
$ws.Range("M74").Value = -3500.5
$ws.Range("H76").Value = 3316.17
$ws.Range("I76").Value = 2948.9773
$ws.Range("K76").Value = 2948.9773
$ws.Range("M76").Value = -2633.9773
$ws.Range("H77").Value = 4585.3
$ws.Range("I77").Value = 4436.5
$ws.Range("K77").Value = 22182.5
$ws.Range("M77").Value = -17502.5
$ws.Range("H79").Value = 3316.17
$ws.Range("I79").Value = 2948.9773
$ws.Range("K79").Value = 2948.9773
$ws.Range("M79").Value = -1856.9773
$ws.Range("H100").Value = 25643660
$ws.Range("I100").Value = 47620996
$ws.Range("K100").Value = 47620996
$ws.Range("M100").Value = -47620455
$ws.Range("H132").Value = 33181.227
$ws.Range("I132").Value = 36668.5
$ws.Range("J132").Value = 633.3333
$ws.Range("K132").Value = 110005.5
$ws.Range("L132").Value = 1899.9999
$ws.Range("M132").Value = -107475.5
$ws.Range("N132").Value = -6959.9999
$ws.Range("H137").Value = 1511.5151
$ws.Range("I137").Value = 1145.9286
$ws.Range("J137").Value = 3558.8
$ws.Range("K137").Value = 3437.7858
$ws.Range("L137").Value = 10676.4
$ws.Range("M137").Value = -887.7857999999997
$ws.Range("N137").Value = -15776.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3269.6584
$ws.Range("I32").Value = 2339.4473
$ws.Range("K32").Value = 2339.4473
$ws.Range("M32").Value = -2052.4473
$ws.Range("H61").Value = 1840.4
$ws.Range("I61").Value = 1831.2632
$ws.Range("K61").Value = 1831.2632
$ws.Range("M61").Value = -1619.2632
$ws.Range("H97").Value = 1604.1034
$ws.Range("I97").Value = 1292.8695
$ws.Range("K97").Value = 1292.8695
$ws.Range("M97").Value = -796.8695
$ws.Range("H132").Value = 2118
$ws.Range("I132").Value = 1232.4
$ws.Range("J132").Value = 3080.6086
$ws.Range("K132").Value = 3697.2
$ws.Range("L132").Value = 9241.825800000001
$ws.Range("M132").Value = -1167.2
$ws.Range("N132").Value = -14301.8258
$ws.Range("H136").Value = 1840.4
$ws.Range("I136").Value = 1831.2632
$ws.Range("K136").Value = 5493.7896
$ws.Range("M136").Value = -2943.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1691.4546
$ws.Range("I86").Value = 1587.375
$ws.Range("J86").Value = 1969
$ws.Range("K86").Value = 1587.375
$ws.Range("L86").Value = 1969
$ws.Range("M86").Value = -464.375
$ws.Range("N86").Value = -4215
$ws.Range("H89").Value = 1691.4546
$ws.Range("I89").Value = 1587.375
$ws.Range("J89").Value = 1969
$ws.Range("K89").Value = 7936.875
$ws.Range("L89").Value = 9845
$ws.Range("M89").Value = -2320.875
$ws.Range("N89").Value = -21077
$ws.Range("H132").Value = 28270
$ws.Range("J132").Value = 28270
$ws.Range("L132").Value = 28270
$ws.Range("N132").Value = -38390
$ws.Range("H134").Value = 2364.5186
$ws.Range("I134").Value = 1307.5294
$ws.Range("K134").Value = 3922.5882
$ws.Range("M134").Value = -1387.5882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2680.1
$ws.Range("I31").Value = 2504.7368
$ws.Range("J31").Value = 2838.762
$ws.Range("K31").Value = 2504.7368
$ws.Range("L31").Value = 2838.762
$ws.Range("M31").Value = -2209.7368
$ws.Range("N31").Value = -3428.762
$ws.Range("H34").Value = 2680.1
$ws.Range("I34").Value = 2504.7368
$ws.Range("J34").Value = 2838.762
$ws.Range("K34").Value = 2504.7368
$ws.Range("L34").Value = 2838.762
$ws.Range("M34").Value = -2302.7368
$ws.Range("N34").Value = -3242.762
$ws.Range("H58").Value = 1687.561
$ws.Range("I58").Value = 1294.88
$ws.Range("J58").Value = 2301.125
$ws.Range("K58").Value = 1294.88
$ws.Range("L58").Value = 2301.125
$ws.Range("M58").Value = -1091.88
$ws.Range("N58").Value = -2707.125
$ws.Range("H62").Value = 4030.3076
$ws.Range("I62").Value = 2319
$ws.Range("K62").Value = 2319
$ws.Range("M62").Value = -1695
$ws.Range("H65").Value = 4030.3076
$ws.Range("I65").Value = 2319
$ws.Range("K65").Value = 11595
$ws.Range("M65").Value = -8475
$ws.Range("H81").Value = 26950
$ws.Range("J81").Value = 26950
$ws.Range("L81").Value = 26950
$ws.Range("N81").Value = -28946
$ws.Range("H84").Value = 26950
$ws.Range("J84").Value = 26950
$ws.Range("L84").Value = 80850
$ws.Range("N84").Value = -90834
$ws.Range("H132").Value = 2648.95
$ws.Range("I132").Value = 1115.6364
$ws.Range("J132").Value = 4523
$ws.Range("K132").Value = 3346.9092
$ws.Range("L132").Value = 13569
$ws.Range("M132").Value = -816.9092000000001
$ws.Range("N132").Value = -18629
$ws.Range("H134").Value = 2762.4546
$ws.Range("I134").Value = 2786.25
$ws.Range("J134").Value = 2699
$ws.Range("K134").Value = 8358.75
$ws.Range("L134").Value = 8097
$ws.Range("M134").Value = -5823.75
$ws.Range("N134").Value = -13167
$ws.Range("H136").Value = 1687.561
$ws.Range("I136").Value = 1294.88
$ws.Range("J136").Value = 2301.125
$ws.Range("K136").Value = 3884.64
$ws.Range("L136").Value = 6903.375
$ws.Range("M136").Value = -1334.64
$ws.Range("N136").Value = -12003.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1923.1
$ws.Range("I122").Value = 1700.175
$ws.Range("K122").Value = 5100.525
$ws.Range("M122").Value = -2650.525
$ws.Range("H126").Value = 1975.0303
$ws.Range("I126").Value = 1601.5238
$ws.Range("K126").Value = 4804.5714
$ws.Range("M126").Value = -2334.5714
$ws.Range("H132").Value = 2268.1516
$ws.Range("I132").Value = 1921.12
$ws.Range("J132").Value = 3352.625
$ws.Range("K132").Value = 5763.36
$ws.Range("L132").Value = 10057.875
$ws.Range("M132").Value = -3233.36
$ws.Range("N132").Value = -15117.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3240.05
$ws.Range("I122").Value = 2099.182
$ws.Range("J122").Value = 4634.4443
$ws.Range("K122").Value = 6297.545999999999
$ws.Range("L122").Value = 13903.3329
$ws.Range("M122").Value = -3847.545999999999
$ws.Range("N122").Value = -18803.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 20458.166
$ws.Range("I76").Value = 3028.5
$ws.Range("J76").Value = 29173
$ws.Range("K76").Value = 3028.5
$ws.Range("L76").Value = 29173
$ws.Range("M76").Value = -2713.5
$ws.Range("N76").Value = -29803
$ws.Range("H79").Value = 20458.166
$ws.Range("I79").Value = 3028.5
$ws.Range("J79").Value = 29173
$ws.Range("K79").Value = 3028.5
$ws.Range("L79").Value = 29173
$ws.Range("M79").Value = -1936.5
$ws.Range("N79").Value = -31357
$ws.Range("H122").Value = 29824.389
$ws.Range("I122").Value = 41729
$ws.Range("J122").Value = 2768.4546
$ws.Range("K122").Value = 125187
$ws.Range("L122").Value = 8305.363799999999
$ws.Range("M122").Value = -122737
$ws.Range("N122").Value = -13205.3638
$ws.Range("H132").Value = 3183.923
$ws.Range("I132").Value = 2961.625
$ws.Range("J132").Value = 3539.6
$ws.Range("K132").Value = 8884.875
$ws.Range("L132").Value = 10618.8
$ws.Range("M132").Value = -6354.875
$ws.Range("N132").Value = -15678.8
$ws.Range("H136").Value = 2979.2
$ws.Range("I136").Value = 1063.2354
$ws.Range("J136").Value = 7050.625
$ws.Range("K136").Value = 3189.7062
$ws.Range("L136").Value = 21151.875
$ws.Range("M136").Value = -639.7062000000001
$ws.Range("N136").Value = -26251.875

